# "Last Commit in Magic"
# Rename Sheet2 -> cred, populate it with a username/password table,
# add a hyperlink on the password cell, resize its columns, move the
# active tab/selection to Sheet2 ("cred"), and update Sheet1's selection.

$wb = $excel.ActiveWorkbook

# --- Sheet2: rename to "cred" and fill with credential data ---
$ws2 = $wb.Worksheets.Item("Sheet2")
$ws2.Name = "cred"

$ws2.Range("A1").Value = "username"
$ws2.Range("B1").Value = "Password"
$ws2.Range("A2").Value = "Admin"
$ws2.Range("B2").Value = 123
$ws2.Range("A3").Value = "Deepakbiet1991"
$ws2.Range("B3").Value = "Kartikey1991@"
$ws2.Range("A4").Value = "Huest"
$ws2.Range("B4").Value = "Guespwd"

# Hyperlink on the password value in row 3, then restore the default
# (non-hyperlink) cell style so the cell keeps its plain formatting.
$ws2.Hyperlinks.Add($ws2.Range("B3"), "http://example.com") | Out-Null
$ws2.Range("B3").Style = "Normal"

# Column widths (closest representable values to the authored widths).
$ws2.Columns.Item(1).ColumnWidth = 16.75
$ws2.Columns.Item(2).ColumnWidth = 11.5

# --- Sheet1: clear tab selection, move cell selection to E6 ---
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws1.Range("E6").Select() | Out-Null

# --- Make Sheet2 ("cred") the active/selected tab with B2 selected ---
$ws2.Activate() | Out-Null
$ws2.Range("B2").Select() | Out-Null
